$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the Job column labels (A2:A7), which were shared-string letters
# A-F, with plain numbers 1-6.
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5
$ws.Range("A7").Value = 6

# Move the selection to A8
$ws.Range("A8").Select()
